$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H17").Value = 2290.1904
$ws.Range("J17").Value = 2290.1904
$ws.Range("L17").Value = 6870.5712
$ws.Range("N17").Value = -7206.5712
$ws.Range("H32").Value = 1016.6667
$ws.Range("J32").Value = 1016.6667
$ws.Range("L32").Value = 1016.6667
$ws.Range("N32").Value = -1668.6667
$ws.Range("H40").Value = 5625.1875
$ws.Range("J40").Value = 3000.7778
$ws.Range("L40").Value = 3000.7778
$ws.Range("N40").Value = -3350.7778
$ws.Range("H51").Value = 3721.3
$ws.Range("I51").Value = 4001
$ws.Range("J51").Value = 3690.2222
$ws.Range("K51").Value = 4001
$ws.Range("L51").Value = 3690.2222
$ws.Range("M51").Value = -3517
$ws.Range("N51").Value = -4658.2222
$ws.Range("H74").Value = 6634.1304
$ws.Range("J74").Value = 7535.1665
$ws.Range("L74").Value = 7535.1665
$ws.Range("N74").Value = -9407.166499999999
$ws.Range("H77").Value = 6634.1304
$ws.Range("J77").Value = 7535.1665
$ws.Range("L77").Value = 37675.8325
$ws.Range("N77").Value = -47035.8325
$ws.Range("H86").Value = 3037.926
$ws.Range("J86").Value = 3572.8
$ws.Range("L86").Value = 3572.8
$ws.Range("N86").Value = -5818.8
$ws.Range("H89").Value = 3037.926
$ws.Range("J89").Value = 3572.8
$ws.Range("L89").Value = 17864
$ws.Range("N89").Value = -29096
$ws.Range("H98").Value = 1130.6666
$ws.Range("I98").Value = 1006.1739
$ws.Range("K98").Value = 1006.1739
$ws.Range("M98").Value = 491.8261
$ws.Range("H107").Value = 532.8570999999999
$ws.Range("I107").Value = 359.54544
$ws.Range("J107").Value = 1168.3334
$ws.Range("K107").Value = 359.54544
$ws.Range("L107").Value = 1168.3334
$ws.Range("M107").Value = 1560.45456
$ws.Range("N107").Value = -5008.3334
$ws.Range("H112").Value = 5338.017
$ws.Range("I112").Value = 1773
$ws.Range("J112").Value = 5897.2354
$ws.Range("K112").Value = 5319
$ws.Range("L112").Value = 17691.7062
$ws.Range("M112").Value = -4211
$ws.Range("N112").Value = -19907.7062
$ws.Range("H116").Value = 307203.72
$ws.Range("I116").Value = 110566.664
$ws.Range("K116").Value = 110566.664
$ws.Range("M116").Value = -107124.664
$ws.Range("H122").Value = 1130.6666
$ws.Range("I122").Value = 1006.1739
$ws.Range("K122").Value = 3018.5217
$ws.Range("M122").Value = -568.5217000000002
$ws.Range("H125").Value = 1989
$ws.Range("I125").Value = 1808.4286
$ws.Range("J125").Value = 2241.8
$ws.Range("K125").Value = 16275.8574
$ws.Range("L125").Value = 20176.2
$ws.Range("M125").Value = -13815.8574
$ws.Range("N125").Value = -25096.2
$ws.Range("H132").Value = 587609.75
$ws.Range("I132").Value = 782147
$ws.Range("K132").Value = 2346441
$ws.Range("M132").Value = -2343911
$ws.Range("H135").Value = 942.3469
$ws.Range("I135").Value = 946.4773
$ws.Range("J135").Value = 906
$ws.Range("K135").Value = 8518.295700000001
$ws.Range("L135").Value = 8154
$ws.Range("M135").Value = -5983.295700000001
$ws.Range("N135").Value = -13224
$ws.Range("H137").Value = 2363.3794
$ws.Range("I137").Value = 1984.5264
$ws.Range("J137").Value = 3083.2
$ws.Range("K137").Value = 5953.5792
$ws.Range("L137").Value = 9249.599999999999
$ws.Range("M137").Value = -3403.5792
$ws.Range("N137").Value = -14349.6
$ws.Range("H138").Value = 1706.1587
$ws.Range("I138").Value = 1134.8837
$ws.Range("K138").Value = 3404.6511
$ws.Range("M138").Value = 1735.3489

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7813968
$ws.Range("J32").Value = 6896.6
$ws.Range("L32").Value = 6896.6
$ws.Range("N32").Value = -7470.6
$ws.Range("H37").Value = 67200
$ws.Range("H45").Value = 2814.0527
$ws.Range("J45").Value = 1467.8334
$ws.Range("L45").Value = 1467.8334
$ws.Range("N45").Value = -2221.8334
$ws.Range("H61").Value = 1518480.5
$ws.Range("I61").Value = 3705554
$ws.Range("J61").Value = 4352.615
$ws.Range("K61").Value = 3705554
$ws.Range("L61").Value = 4352.615
$ws.Range("M61").Value = -3705342
$ws.Range("N61").Value = -4776.615
$ws.Range("H74").Value = 2018939.4
$ws.Range("I74").Value = 2406267.8
$ws.Range("J74").Value = 4831.6
$ws.Range("K74").Value = 2406267.8
$ws.Range("L74").Value = 4831.6
$ws.Range("M74").Value = -2405393.8
$ws.Range("N74").Value = -6579.6
$ws.Range("H77").Value = 2018939.4
$ws.Range("I77").Value = 2406267.8
$ws.Range("J77").Value = 4831.6
$ws.Range("K77").Value = 12031339
$ws.Range("L77").Value = 24158
$ws.Range("M77").Value = -12026971
$ws.Range("N77").Value = -32894
$ws.Range("H88").Value = 1190.8
$ws.Range("I88").Value = 1312.25
$ws.Range("J88").Value = 1146.6364
$ws.Range("K88").Value = 1312.25
$ws.Range("L88").Value = 1146.6364
$ws.Range("M88").Value = -906.25
$ws.Range("N88").Value = -1958.6364
$ws.Range("H91").Value = 1190.8
$ws.Range("I91").Value = 1312.25
$ws.Range("J91").Value = 1146.6364
$ws.Range("K91").Value = 1312.25
$ws.Range("L91").Value = 1146.6364
$ws.Range("M91").Value = 91.75
$ws.Range("N91").Value = -3954.6364
$ws.Range("H122").Value = 3154.862
$ws.Range("I122").Value = 3000.2727
$ws.Range("J122").Value = 3640.7144
$ws.Range("K122").Value = 9000.8181
$ws.Range("L122").Value = 10922.1432
$ws.Range("M122").Value = -6550.8181
$ws.Range("N122").Value = -15822.1432
$ws.Range("H132").Value = 834641.2
$ws.Range("I132").Value = 834641.2
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2503923.6
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2501393.6
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 1518480.5
$ws.Range("I136").Value = 3705554
$ws.Range("J136").Value = 4352.615
$ws.Range("K136").Value = 11116662
$ws.Range("L136").Value = 13057.845
$ws.Range("M136").Value = -11114112
$ws.Range("N136").Value = -18157.845
$ws.Range("H138").Value = 199999
$ws.Range("J138").Value = 199999
$ws.Range("L138").Value = 199999
$ws.Range("N138").Value = -210279

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 827.2222
$ws.Range("J7").Value = 1166.6666
$ws.Range("L7").Value = 1166.6666
$ws.Range("N7").Value = -1392.6666
$ws.Range("H20").Value = 1698.2222
$ws.Range("I20").Value = 1662.8125
$ws.Range("K20").Value = 1662.8125
$ws.Range("M20").Value = -1415.8125
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("H86").Value = 3799.6667
$ws.Range("I86").Value = 3500
$ws.Range("K86").Value = 3500
$ws.Range("M86").Value = -2377
$ws.Range("H89").Value = 3799.6667
$ws.Range("I89").Value = 3500
$ws.Range("K89").Value = 17500
$ws.Range("M89").Value = -11884
$ws.Range("H99").Value = 77502.625
$ws.Range("I99").Value = 86001.664
$ws.Range("K99").Value = 86001.664
$ws.Range("M99").Value = -84503.664
$ws.Range("H127").Value = 45000
$ws.Range("J127").Value = 45000
$ws.Range("L127").Value = 45000
$ws.Range("N127").Value = -54920
$ws.Range("H129").Value = 74000
$ws.Range("J129").Value = 74000
$ws.Range("L129").Value = 74000
$ws.Range("N129").Value = -84000
$ws.Range("H130").Value = 69566.664
$ws.Range("J130").Value = 69566.664
$ws.Range("L130").Value = 69566.664
$ws.Range("N130").Value = -79606.664
$ws.Range("H134").Value = 995508
$ws.Range("I134").Value = 1702709.2
$ws.Range("J134").Value = 5426.3
$ws.Range("K134").Value = 5108127.6
$ws.Range("L134").Value = 16278.9
$ws.Range("M134").Value = -5105592.6
$ws.Range("N134").Value = -21348.9

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("H31").Value = 93283.92
$ws.Range("I31").Value = 147928.6
$ws.Range("J31").Value = 26496
$ws.Range("K31").Value = 147928.6
$ws.Range("L31").Value = 26496
$ws.Range("M31").Value = -147633.6
$ws.Range("N31").Value = -27086
$ws.Range("H34").Value = 93283.92
$ws.Range("I34").Value = 147928.6
$ws.Range("J34").Value = 26496
$ws.Range("K34").Value = 147928.6
$ws.Range("L34").Value = 26496
$ws.Range("M34").Value = -147726.6
$ws.Range("N34").Value = -26900
$ws.Range("H58").Value = 1126550.6
$ws.Range("J58").Value = 5539.125
$ws.Range("L58").Value = 5539.125
$ws.Range("N58").Value = -5945.125
$ws.Range("H62").Value = 7148.125
$ws.Range("J62").Value = 3899.5
$ws.Range("L62").Value = 3899.5
$ws.Range("N62").Value = -5147.5
$ws.Range("H65").Value = 7148.125
$ws.Range("J65").Value = 3899.5
$ws.Range("L65").Value = 19497.5
$ws.Range("N65").Value = -25737.5
$ws.Range("H80").Value = 47995.5
$ws.Range("J80").Value = 47995.5
$ws.Range("L80").Value = 47995.5
$ws.Range("N80").Value = -50241.5
$ws.Range("H83").Value = 47995.5
$ws.Range("J83").Value = 47995.5
$ws.Range("L83").Value = 143986.5
$ws.Range("N83").Value = -155218.5
$ws.Range("H87").Value = 123553.336
$ws.Range("J87").Value = 123553.336
$ws.Range("L87").Value = 123553.336
$ws.Range("N87").Value = -125925.336
$ws.Range("H88").Value = 16315.333
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 16315.333
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 16315.333
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -17127.333
$ws.Range("H90").Value = 123553.336
$ws.Range("J90").Value = 123553.336
$ws.Range("L90").Value = 370660.008
$ws.Range("N90").Value = -382516.008
$ws.Range("H91").Value = 16315.333
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 16315.333
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 16315.333
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -19123.333
$ws.Range("H99").Value = 4253.5713
$ws.Range("I99").Value = 4561.3335
$ws.Range("J99").Value = 4022.75
$ws.Range("K99").Value = 4561.3335
$ws.Range("L99").Value = 4022.75
$ws.Range("M99").Value = -3063.3335
$ws.Range("N99").Value = -7018.75
$ws.Range("H105").Value = 17529.273
$ws.Range("I105").Value = 28212.691
$ws.Range("K105").Value = 28212.691
$ws.Range("M105").Value = -26465.691
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H122").Value = 1754.9565
$ws.Range("I122").Value = 1234.2106
$ws.Range("J122").Value = 4228.5
$ws.Range("K122").Value = 3702.6318
$ws.Range("L122").Value = 12685.5
$ws.Range("M122").Value = -1252.6318
$ws.Range("N122").Value = -17585.5
$ws.Range("H126").Value = 4253.5713
$ws.Range("I126").Value = 4561.3335
$ws.Range("J126").Value = 4022.75
$ws.Range("K126").Value = 13684.0005
$ws.Range("L126").Value = 12068.25
$ws.Range("M126").Value = -11214.0005
$ws.Range("N126").Value = -17008.25
$ws.Range("H132").Value = 8077774
$ws.Range("I132").Value = 14889.296
$ws.Range("J132").Value = 62502250
$ws.Range("K132").Value = 44667.888
$ws.Range("L132").Value = 187506750
$ws.Range("M132").Value = -42137.888
$ws.Range("N132").Value = -187511810
$ws.Range("H134").Value = 17824.904
$ws.Range("I134").Value = 28779.334
$ws.Range("J134").Value = 3219
$ws.Range("K134").Value = 86338.00199999999
$ws.Range("L134").Value = 9657
$ws.Range("M134").Value = -83803.00199999999
$ws.Range("N134").Value = -14727
$ws.Range("H136").Value = 1126550.6
$ws.Range("J136").Value = 5539.125
$ws.Range("L136").Value = 16617.375
$ws.Range("N136").Value = -21717.375

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 18718.857
$ws.Range("I87").Value = 13996
$ws.Range("K87").Value = 41988
$ws.Range("M87").Value = -40740
$ws.Range("H90").Value = 18718.857
$ws.Range("I90").Value = 13996
$ws.Range("K90").Value = 125964
$ws.Range("M90").Value = -119724
$ws.Range("H136").Value = 4844.7
$ws.Range("I136").Value = 4844.7
$ws.Range("K136").Value = 14534.1
$ws.Range("M136").Value = -9434.099999999999
$ws.Range("H139").Value = 848.35297
$ws.Range("I139").Value = 848.35297
$ws.Range("K139").Value = 2545.05891
$ws.Range("M139").Value = 2594.94109

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 143.66667
$ws.Range("I2").Value = 55.64706
$ws.Range("K2").Value = 55.64706
$ws.Range("M2").Value = 57.35294
$ws.Range("H43").Value = 7983.1816
$ws.Range("I43").Value = 3477
$ws.Range("K43").Value = 3477
$ws.Range("M43").Value = -3326
$ws.Range("H70").Value = 7728
$ws.Range("I70").Value = 8507.714
$ws.Range("K70").Value = 8507.714
$ws.Range("M70").Value = -8237.714
$ws.Range("H73").Value = 7728
$ws.Range("I73").Value = 8507.714
$ws.Range("K73").Value = 8507.714
$ws.Range("M73").Value = -7571.714
$ws.Range("H80").Value = 5490.1665
$ws.Range("J80").Value = 6485.4287
$ws.Range("L80").Value = 6485.4287
$ws.Range("N80").Value = -8481.4287
$ws.Range("H83").Value = 5490.1665
$ws.Range("J83").Value = 6485.4287
$ws.Range("L83").Value = 32427.1435
$ws.Range("N83").Value = -42411.14350000001
$ws.Range("H97").Value = 1499.3
$ws.Range("I97").Value = 729.9655
$ws.Range("J97").Value = 3527.5454
$ws.Range("K97").Value = 729.9655
$ws.Range("L97").Value = 3527.5454
$ws.Range("M97").Value = -233.9655
$ws.Range("N97").Value = -4519.5454
$ws.Range("H113").Value = 3067.2646
$ws.Range("J113").Value = 3831.9375
$ws.Range("L113").Value = 3831.9375
$ws.Range("N113").Value = -8171.9375
$ws.Range("H122").Value = 27783760
$ws.Range("I122").Value = 35716620
$ws.Range("J122").Value = 18752
$ws.Range("K122").Value = 107149860
$ws.Range("L122").Value = 56256
$ws.Range("M122").Value = -107147410
$ws.Range("N122").Value = -61156
$ws.Range("H132").Value = 1508704.5
$ws.Range("I132").Value = 2010198.6
$ws.Range("J132").Value = 4222
$ws.Range("K132").Value = 6030595.800000001
$ws.Range("L132").Value = 12666
$ws.Range("M132").Value = -6028065.800000001
$ws.Range("N132").Value = -17726

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2667.3333
$ws.Range("I7").Value = 2400.8
$ws.Range("K7").Value = 2400.8
$ws.Range("M7").Value = -2288.8
$ws.Range("H16").Value = 7864
$ws.Range("I16").Value = 9008.5
$ws.Range("J16").Value = 997
$ws.Range("K16").Value = 9008.5
$ws.Range("L16").Value = 997
$ws.Range("M16").Value = -8838.5
$ws.Range("N16").Value = -1337
$ws.Range("H40").Value = 3155.52
$ws.Range("I40").Value = 2599.15
$ws.Range("K40").Value = 2599.15
$ws.Range("M40").Value = -2463.15
$ws.Range("H46").Value = 1062.75
$ws.Range("I46").Value = 1111
$ws.Range("J46").Value = 950.1667
$ws.Range("K46").Value = 1111
$ws.Range("L46").Value = 950.1667
$ws.Range("M46").Value = -923
$ws.Range("N46").Value = -1326.1667
$ws.Range("H55").Value = 933.625
$ws.Range("I55").Value = 277.2857
$ws.Range("K55").Value = 277.2857
$ws.Range("M55").Value = -104.2857
$ws.Range("H68").Value = 4174.1333
$ws.Range("I68").Value = 2085.7144
$ws.Range("K68").Value = 2085.7144
$ws.Range("M68").Value = -1336.7144
$ws.Range("H71").Value = 4174.1333
$ws.Range("I71").Value = 2085.7144
$ws.Range("K71").Value = 10428.572
$ws.Range("M71").Value = -6684.572
$ws.Range("H100").Value = 6314.1665
$ws.Range("I100").Value = 1949.5
$ws.Range("J100").Value = 12424.7
$ws.Range("K100").Value = 1949.5
$ws.Range("L100").Value = 12424.7
$ws.Range("M100").Value = -1408.5
$ws.Range("N100").Value = -13506.7
$ws.Range("H122").Value = 2578.975
$ws.Range("I122").Value = 2552.0908
$ws.Range("K122").Value = 7656.2724
$ws.Range("M122").Value = -5206.2724
$ws.Range("H126").Value = 2667.3333
$ws.Range("I126").Value = 2400.8
$ws.Range("K126").Value = 7202.400000000001
$ws.Range("M126").Value = -4732.400000000001
$ws.Range("H132").Value = 990394
$ws.Range("I132").Value = 1082915.4
$ws.Range("K132").Value = 3248746.2
$ws.Range("M132").Value = -3246216.2
$ws.Range("H136").Value = 84583.25
$ws.Range("I136").Value = 7259.25
$ws.Range("J136").Value = 316555.25
$ws.Range("K136").Value = 21777.75
$ws.Range("L136").Value = 949665.75
$ws.Range("M136").Value = -19227.75
$ws.Range("N136").Value = -954765.75
$ws.Range("H139").Value = 84619
$ws.Range("J139").Value = 93542.8
$ws.Range("L139").Value = 93542.8
$ws.Range("N139").Value = -103822.8

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4791.6665
$ws.Range("J62").Value = 6400
$ws.Range("L62").Value = 6400
$ws.Range("N62").Value = -7648
$ws.Range("H65").Value = 4791.6665
$ws.Range("J65").Value = 6400
$ws.Range("L65").Value = 32000
$ws.Range("N65").Value = -38240
$ws.Range("H81").Value = 1928.9166
$ws.Range("I81").Value = 1897
$ws.Range("K81").Value = 3794
$ws.Range("M81").Value = -2733
$ws.Range("H84").Value = 1928.9166
$ws.Range("I84").Value = 1897
$ws.Range("K84").Value = 18970
$ws.Range("M84").Value = -13666
$ws.Range("H107").Value = 3183.6428
$ws.Range("I107").Value = 1847.05
$ws.Range("K107").Value = 5541.15
$ws.Range("M107").Value = -3621.15
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("I113").Value = 825.8946999999999
$ws.Range("J113").Value = 3104.077
$ws.Range("K113").Value = 2477.6841
$ws.Range("L113").Value = 9312.231
$ws.Range("M113").Value = -307.6840999999999
$ws.Range("N113").Value = -13652.231
$ws.Range("H122").Value = 1892.0883
$ws.Range("I122").Value = 1688.1936
$ws.Range("K122").Value = 5064.5808
$ws.Range("M122").Value = -2614.5808
$ws.Range("H126").Value = 2145.4666
$ws.Range("J126").Value = 3098.5
$ws.Range("L126").Value = 9295.5
$ws.Range("N126").Value = -14235.5
$ws.Range("H131").Value = 66714.14
$ws.Range("J131").Value = 66714.14
$ws.Range("L131").Value = 66714.14
$ws.Range("N131").Value = -76794.14
$ws.Range("H132").Value = 9589377
$ws.Range("I132").Value = 16773745
$ws.Range("J132").Value = 10220.667
$ws.Range("K132").Value = 50321235
$ws.Range("L132").Value = 30662.001
$ws.Range("M132").Value = -50318705
$ws.Range("N132").Value = -35722.001
$ws.Range("H136").Value = 5763597
$ws.Range("I136").Value = 6786288
$ws.Range("K136").Value = 20358864
$ws.Range("M136").Value = -20356314
